$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I31").Value = 499
$ws.Range("K31").Value = 1497
$ws.Range("M31").Value = -1267
$ws.Range("H31").Value = 499
$ws.Range("H48").Value = 1700
$ws.Range("J48").Value = 1800
$ws.Range("N48").Value = -5984
$ws.Range("L48").Value = 5400
$ws.Range("H56").Value = 1700
$ws.Range("L56").Value = 5400
$ws.Range("J56").Value = 1800
$ws.Range("N56").Value = -6468
$ws.Range("K62").Value = 9265010
$ws.Range("I62").Value = 9265010
$ws.Range("M62").Value = -9264386
$ws.Range("H62").Value = 6672760.5
$ws.Range("J62").Value = 6976.2856
$ws.Range("N62").Value = -8224.285599999999
$ws.Range("L62").Value = 6976.2856
$ws.Range("K65").Value = 46325050
$ws.Range("M65").Value = -46321930
$ws.Range("I65").Value = 9265010
$ws.Range("N65").Value = -41121.428
$ws.Range("L65").Value = 34881.428
$ws.Range("H65").Value = 6672760.5
$ws.Range("J65").Value = 6976.2856
$ws.Range("J132").Value = 4378.6
$ws.Range("L132").Value = 13135.8
$ws.Range("H132").Value = 3468.9143
$ws.Range("N132").Value = -18195.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1678.3334
$ws.Range("M2").Value = -1340.3889
$ws.Range("I2").Value = 1453.3889
$ws.Range("K2").Value = 1453.3889
$ws.Range("H5").Value = 2379.6
$ws.Range("K5").Value = 474.5
$ws.Range("M5").Value = -362.5
$ws.Range("I5").Value = 474.5
$ws.Range("L27").Value = 9997
$ws.Range("H27").Value = 9997
$ws.Range("N27").Value = -10365
$ws.Range("J27").Value = 9997
$ws.Range("I33").Value = 18500
$ws.Range("K33").Value = 18500
$ws.Range("M33").Value = -18171
$ws.Range("H33").Value = 18500
$ws.Range("N45").Value = -3631.4
$ws.Range("L45").Value = 2877.4
$ws.Range("J45").Value = 2877.4
$ws.Range("H45").Value = 2714.2307
$ws.Range("K116").Value = 1453.3889
$ws.Range("H116").Value = 1678.3334
$ws.Range("I116").Value = 1453.3889
$ws.Range("M116").Value = 840.6111000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K3").Value = 1453.3889
$ws.Range("I3").Value = 1453.3889
$ws.Range("H3").Value = 1678.3334
$ws.Range("M3").Value = -1339.3889
$ws.Range("M4").Value = -359.5
$ws.Range("K4").Value = 474.5
$ws.Range("I4").Value = 474.5
$ws.Range("H4").Value = 2379.6
$ws.Range("H81").Value = 41395.43
$ws.Range("L81").Value = 39961.332
$ws.Range("N81").Value = -42083.332
$ws.Range("J81").Value = 39961.332
$ws.Range("L84").Value = 119883.996
$ws.Range("H84").Value = 41395.43
$ws.Range("N84").Value = -130491.996
$ws.Range("J84").Value = 39961.332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I58").Value = 3616.3333
$ws.Range("H58").Value = 3676.8333
$ws.Range("J58").Value = 3797.8333
$ws.Range("M58").Value = -3413.3333
$ws.Range("L58").Value = 3797.8333
$ws.Range("N58").Value = -4203.8333
$ws.Range("K58").Value = 3616.3333
$ws.Range("H105").Value = 5945
$ws.Range("I105").Value = 1995
$ws.Range("M105").Value = -248
$ws.Range("K105").Value = 1995
$ws.Range("H134").Value = 3280.0715
$ws.Range("K134").Value = 9822
$ws.Range("M134").Value = -7287
$ws.Range("I134").Value = 3274
$ws.Range("K136").Value = 10848.9999
$ws.Range("J136").Value = 3797.8333
$ws.Range("N136").Value = -16493.4999
$ws.Range("H136").Value = 3676.8333
$ws.Range("M136").Value = -8298.999899999999
$ws.Range("I136").Value = 3616.3333
$ws.Range("L136").Value = 11393.4999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N12").Value = -779.5
$ws.Range("K12").Value = 53.625
$ws.Range("M12").Value = 119.375
$ws.Range("L12").Value = 433.5
$ws.Range("I12").Value = 17.875
$ws.Range("H12").Value = 88.22221999999999
$ws.Range("J12").Value = 144.5
$ws.Range("I23").Value = 114.666664
$ws.Range("H23").Value = 955.26666
$ws.Range("J23").Value = 1515.6666
$ws.Range("K23").Value = 343.999992
$ws.Range("N23").Value = -5016.9998
$ws.Range("M23").Value = -108.999992
$ws.Range("L23").Value = 4546.9998
$ws.Range("M28").Value = -3867.0002
$ws.Range("K28").Value = 4099.0002
$ws.Range("I28").Value = 1366.3334
$ws.Range("H28").Value = 1366.3334
$ws.Range("M40").Value = -651
$ws.Range("H40").Value = 316.66666
$ws.Range("K40").Value = 720
$ws.Range("I40").Value = 180
$ws.Range("I137").Value = 2335.8333
$ws.Range("M137").Value = -1907.499899999999
$ws.Range("K137").Value = 7007.499899999999
$ws.Range("H137").Value = 3213.8333
$ws.Range("I140").Value = 12500931
$ws.Range("H140").Value = 6946796.5
$ws.Range("K140").Value = 37502793
$ws.Range("M140").Value = -37497613
$ws.Range("M141").Value = -4621.75
$ws.Range("I141").Value = 3267.25
$ws.Range("H141").Value = 3713.7
$ws.Range("K141").Value = 9801.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("L97").Value = 7000
$ws.Range("H97").Value = 2764.0435
$ws.Range("I97").Value = 2128.65
$ws.Range("N97").Value = -7992
$ws.Range("J97").Value = 7000
$ws.Range("K97").Value = 2128.65
$ws.Range("M97").Value = -1632.65
$ws.Range("M122").Value = -13346.845
$ws.Range("H122").Value = 5767.6
$ws.Range("I122").Value = 5265.615
$ws.Range("K122").Value = 15796.845
$ws.Range("M132").Value = -87730.45999999999
$ws.Range("I132").Value = 30086.82
$ws.Range("K132").Value = 90260.45999999999
$ws.Range("H132").Value = 20664.594

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value = 4918.8335
$ws.Range("M7").Value = -4825.3687
$ws.Range("L7").Value = 4918.8335
$ws.Range("H7").Value = 4932.92
$ws.Range("N7").Value = -5142.8335
$ws.Range("K7").Value = 4937.3687
$ws.Range("I7").Value = 4937.3687
$ws.Range("H16").Value = 1614.919
$ws.Range("I16").Value = 1582.5769
$ws.Range("K16").Value = 1582.5769
$ws.Range("M16").Value = -1412.5769
$ws.Range("M40").Value = -17481.5
$ws.Range("J40").Value = 0
$ws.Range("H40").Value = 17617.5
$ws.Range("K40").Value = 17617.5
$ws.Range("I40").Value = 17617.5
$ws.Range("N40").ClearContents()
$ws.Range("L40").Value = 0
$ws.Range("H55").Value = 7105.625
$ws.Range("N55").Value = -38596.25
$ws.Range("L55").Value = 38250.25
$ws.Range("J55").Value = 38250.25
$ws.Range("H61").Value = 2149.1177
$ws.Range("K61").Value = 2189.6875
$ws.Range("M61").Value = -1987.6875
$ws.Range("I61").Value = 2189.6875
$ws.Range("K68").Value = 10612.533
$ws.Range("H68").Value = 10074.25
$ws.Range("M68").Value = -9863.532999999999
$ws.Range("I68").Value = 10612.533
$ws.Range("H71").Value = 10074.25
$ws.Range("M71").Value = -49318.66499999999
$ws.Range("I71").Value = 10612.533
$ws.Range("K71").Value = 53062.66499999999
$ws.Range("I93").Value = 2349.9092
$ws.Range("M93").Value = -1101.9092
$ws.Range("K93").Value = 2349.9092
$ws.Range("H93").Value = 66669412
$ws.Range("J93").Value = 250003840
$ws.Range("N93").Value = -250006336
$ws.Range("L93").Value = 250003840
$ws.Range("M113").Value = -19.6875
$ws.Range("H113").Value = 2149.1177
$ws.Range("I113").Value = 2189.6875
$ws.Range("K113").Value = 2189.6875
$ws.Range("L126").Value = 14756.5005
$ws.Range("K126").Value = 14812.1061
$ws.Range("N126").Value = -19696.5005
$ws.Range("I126").Value = 4937.3687
$ws.Range("H126").Value = 4932.92
$ws.Range("M126").Value = -12342.1061
$ws.Range("J126").Value = 4918.8335
$ws.Range("K136").Value = 2923.94109
$ws.Range("J136").Value = 3630.3
$ws.Range("N136").Value = -15990.9
$ws.Range("H136").Value = 1958.2222
$ws.Range("M136").Value = -373.9410899999998
$ws.Range("I136").Value = 974.64703
$ws.Range("L136").Value = 10890.9

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8964
$ws.Range("J62").Value = 11263.167
$ws.Range("N62").Value = -12511.167
$ws.Range("L62").Value = 11263.167
$ws.Range("N65").Value = -62555.835
$ws.Range("L65").Value = 56315.835
$ws.Range("H65").Value = 8964
$ws.Range("J65").Value = 11263.167
$ws.Range("K81").Value = 7851.857
$ws.Range("H81").Value = 8174.294
$ws.Range("M81").Value = -6790.857
$ws.Range("I81").Value = 3925.9285
$ws.Range("I84").Value = 3925.9285
$ws.Range("M84").Value = -33955.285
$ws.Range("K84").Value = 39259.285
$ws.Range("H84").Value = 8174.294
$ws.Range("H96").Value = 2421.5483
$ws.Range("M96").Value = -1391.889
$ws.Range("I96").Value = 2764.889
$ws.Range("K96").Value = 2764.889
$ws.Range("N107").Value = -4737
$ws.Range("J107").Value = 299
$ws.Range("H107").Value = 659.4
$ws.Range("L107").Value = 897
$ws.Range("M122").Value = -6409.2145
$ws.Range("H122").Value = 3022.8667
$ws.Range("I122").Value = 2953.0715
$ws.Range("K122").Value = 8859.2145
$ws.Range("K126").Value = 11872.6671
$ws.Range("I126").Value = 3957.5557
$ws.Range("H126").Value = 4329.2144
$ws.Range("M126").Value = -9402.667099999999

Write-Output "Applied all market-price updates."